$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("E2").Value = 22
$ws.Range("G2").Value = [double]"2.753020034163001e-12"
$ws.Range("H2").Value = [double]"8.031900426556314e-12"
$ws.Range("I2").Value = 0.1391838349699088
$ws.Range("K2").Value = 38.7665550046719
$ws.Range("L2").Value = "[27.787891635915592, 49.74521837342821]"
$ws.Range("M2").Value = [double]"5.273825820495404e-11"
$ws.Range("N2").Value = [double]"5.273825820495404e-11"
$ws.Range("O2").Value = 1.50318447288881
$ws.Range("P2").Value = "[1.1635528346628856, 1.8428161111147343]"
$ws.Range("Q2").Value = [double]"1.332267629550188e-15"
$ws.Range("R2").Value = [double]"1.332267629550188e-15"
$ws.Range("S2").Value = 54.65485085037316
$ws.Range("T2").Value = "[47.71791256022658, 61.59178914051974]"
$ws.Range("W2").Value = 16.73673673673674
$ws.Range("X2").Value = 15.54754754754755
$ws.Range("Y2").Value = 17.92592592592593

# Row 3 updates
$ws.Range("E3").Value = 24.69000000000042
$ws.Range("H3").Value = [double]"3.51058663913093e-16"
$ws.Range("K3").Value = 49.89521492954081
$ws.Range("L3").Value = "[40.124360729582556, 59.66606912949906]"
$ws.Range("O3").Value = -2.566105711040311
$ws.Range("P3").Value = "[-2.76736890406308, -2.364842518017543]"
$ws.Range("S3").Value = 63.94937950423303
$ws.Range("T3").Value = "[58.88155936327668, 69.01719964518938]"
$ws.Range("W3").Value = 10.08360360360378
$ws.Range("X3").Value = 9.292732732732896
$ws.Range("Y3").Value = 10.87447447447465
